$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix B12: it was stored as text "1000271912"; it should be a real number ---
$ws.Cells.Item(12, 2).Value = 1000271912

# --- 2) Add new row 13 with the negotiation log entry dated 2025-12-15 10:39:04 ---
# Start from a copy of row 12 so the new row inherits the same (unstyled) text
# formatting/types used throughout the sheet (plain text cells, no direct formatting).
$ws.Rows.Item(12).Copy()
$ws.Rows.Item(13).PasteSpecial()

$ws.Cells.Item(13, 1).Value  = "2025-12-15 10:39:04"
$ws.Cells.Item(13, 3).Value  = "Leidy"
$ws.Cells.Item(13, 4).Value  = "TARJETA DE CRÉDITO"
$ws.Cells.Item(13, 5).Value  = "****9053"
$ws.Cells.Item(13, 6).Value  = "PRORROGA CON PAGO"
$ws.Cells.Item(13, 7).Value  = "12 cuotas"
$ws.Cells.Item(13, 8).Value  = "35.230.127.150"
$ws.Cells.Item(13, 9).Value  = "The Dalles"
$ws.Cells.Item(13, 10).Value = "Oregon"
$ws.Cells.Item(13, 11).Value = "United States"
$ws.Cells.Item(13, 12).Value = "2025-12-15 10:39:04"
$ws.Cells.Item(13, 13).Value = "*****9053"
$ws.Cells.Item(13, 14).Value = "35.230.127.150"

# Cedula (column B) is the text "1000135120" rather than a number, so force text
# storage before assigning to avoid automatic numeric conversion, then drop the
# number-format override so no extra direct formatting remains on the cell.
$ws.Cells.Item(13, 2).NumberFormat = "@"
$ws.Cells.Item(13, 2).Value = "1000135120"
$ws.Cells.Item(13, 2).Style = "Normal"
